$wb = $excel.ActiveWorkbook

# This script applies updated market-board price/profit figures
# (columns H-N) across multiple job sheets, as produced by the
# scheduled market-data refresh runner.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 6333.3335
$ws.Range("I12").Value = 5000
$ws.Range("J12").Value = 7000
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = -4830
$ws.Range("N12").Value = -7340
$ws.Range("H98").Value = 879.4286
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 1000
$ws.Range("N98").Value = -3996
$ws.Range("H122").Value = 879.4286
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 8800290
$ws.Range("I132").Value = 10102461
$ws.Range("K132").Value = 30307383
$ws.Range("M132").Value = -30304853
$ws.Range("H137").Value = 5989.8184
$ws.Range("I137").Value = 1353.375
$ws.Range("K137").Value = 4060.125
$ws.Range("M137").Value = -1510.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27029934
$ws.Range("I2").Value = 35716984
$ws.Range("K2").Value = 35716984
$ws.Range("M2").Value = -35716871
$ws.Range("H32").Value = 3303
$ws.Range("I32").Value = 1696.0864
$ws.Range("K32").Value = 1696.0864
$ws.Range("M32").Value = -1409.0864
$ws.Range("H74").Value = 33689.863
$ws.Range("I74").Value = 54223.74
$ws.Range("J74").Value = 12015.223
$ws.Range("K74").Value = 54223.74
$ws.Range("L74").Value = 12015.223
$ws.Range("M74").Value = -53349.74
$ws.Range("N74").Value = -13763.223
$ws.Range("H77").Value = 33689.863
$ws.Range("I77").Value = 54223.74
$ws.Range("J77").Value = 12015.223
$ws.Range("K77").Value = 271118.7
$ws.Range("L77").Value = 60076.115
$ws.Range("M77").Value = -266750.7
$ws.Range("N77").Value = -68812.11499999999
$ws.Range("H116").Value = 27029934
$ws.Range("I116").Value = 35716984
$ws.Range("K116").Value = 35716984
$ws.Range("M116").Value = -35714690
$ws.Range("H135").Value = 45406
$ws.Range("J135").Value = 45406
$ws.Range("L135").Value = 45406
$ws.Range("N135").Value = -55546
$ws.Range("H139").Value = 87763
$ws.Range("J139").Value = 87763
$ws.Range("L139").Value = 87763
$ws.Range("N139").Value = -98043

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27029934
$ws.Range("I3").Value = 35716984
$ws.Range("K3").Value = 35716984
$ws.Range("M3").Value = -35716870
$ws.Range("H86").Value = 2153.8462
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2153.8462
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H94").Value = 4399.75
$ws.Range("I94").Value = 4466.3335
$ws.Range("K94").Value = 4466.3335
$ws.Range("M94").Value = -4015.3335
$ws.Range("H99").Value = 2218.8235
$ws.Range("I99").Value = 1228.3334
$ws.Range("J99").Value = 4596
$ws.Range("K99").Value = 1228.3334
$ws.Range("L99").Value = 4596
$ws.Range("M99").Value = 269.6666
$ws.Range("N99").Value = -7592
$ws.Range("H132").Value = 117494.5
$ws.Range("J132").Value = 117494.5
$ws.Range("L132").Value = 117494.5
$ws.Range("N132").Value = -127614.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25417.156
$ws.Range("I31").Value = 48438
$ws.Range("K31").Value = 48438
$ws.Range("M31").Value = -48143
$ws.Range("H34").Value = 25417.156
$ws.Range("I34").Value = 48438
$ws.Range("K34").Value = 48438
$ws.Range("M34").Value = -48236
$ws.Range("H58").Value = 2539.8696
$ws.Range("I58").Value = 2393.513
$ws.Range("K58").Value = 2393.513
$ws.Range("M58").Value = -2190.513
$ws.Range("H86").Value = 10288.652
$ws.Range("I86").Value = 9241.546
$ws.Range("J86").Value = 11248.5
$ws.Range("K86").Value = 9241.546
$ws.Range("L86").Value = 11248.5
$ws.Range("M86").Value = -8118.546
$ws.Range("N86").Value = -13494.5
$ws.Range("H89").Value = 10288.652
$ws.Range("I89").Value = 9241.546
$ws.Range("J89").Value = 11248.5
$ws.Range("K89").Value = 46207.73
$ws.Range("L89").Value = 56242.5
$ws.Range("M89").Value = -40591.73
$ws.Range("N89").Value = -67474.5
$ws.Range("H136").Value = 2539.8696
$ws.Range("I136").Value = 2393.513
$ws.Range("K136").Value = 7180.539
$ws.Range("M136").Value = -4630.539

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 446.60526
$ws.Range("I6").Value = 451.9189
$ws.Range("K6").Value = 1355.7567
$ws.Range("M6").Value = -1242.7567
$ws.Range("H10").Value = 21.181818
$ws.Range("I10").Value = 21.444445
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 64.33333500000001
$ws.Range("L10").Value = 60
$ws.Range("M10").Value = 74.66666499999999
$ws.Range("N10").Value = -338
$ws.Range("H15").Value = 287.1111
$ws.Range("I15").Value = 108.5
$ws.Range("J15").Value = 430
$ws.Range("K15").Value = 325.5
$ws.Range("L15").Value = 1290
$ws.Range("M15").Value = -185.5
$ws.Range("N15").Value = -1570
$ws.Range("H44").Value = 850.5
$ws.Range("I44").Value = 520.6
$ws.Range("K44").Value = 1561.8
$ws.Range("M44").Value = -1163.8
$ws.Range("H102").Value = 5972.222
$ws.Range("J102").Value = 7692.3076
$ws.Range("L102").Value = 23076.9228
$ws.Range("N102").Value = -27944.9228
$ws.Range("H131").Value = 24974.273
$ws.Range("I131").Value = 72320.57000000001
$ws.Range("J131").Value = 2879.3333
$ws.Range("K131").Value = 216961.71
$ws.Range("L131").Value = 8637.999899999999
$ws.Range("M131").Value = -211921.71
$ws.Range("N131").Value = -18717.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 209.89131
$ws.Range("I2").Value = 146.6
$ws.Range("J2").Value = 328.5625
$ws.Range("K2").Value = 146.6
$ws.Range("L2").Value = 328.5625
$ws.Range("M2").Value = -33.59999999999999
$ws.Range("N2").Value = -554.5625
$ws.Range("H57").Value = 10105.1
$ws.Range("J57").Value = 14999
$ws.Range("L57").Value = 14999
$ws.Range("N57").Value = -16639
$ws.Range("H122").Value = 3572
$ws.Range("I122").Value = 3015.6
$ws.Range("J122").Value = 3969.4285
$ws.Range("K122").Value = 9046.799999999999
$ws.Range("L122").Value = 11908.2855
$ws.Range("M122").Value = -6596.799999999999
$ws.Range("N122").Value = -16808.2855
$ws.Range("H132").Value = 3072.6086
$ws.Range("I132").Value = 3067.875
$ws.Range("K132").Value = 9203.625
$ws.Range("M132").Value = -6673.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 26318146
$ws.Range("I16").Value = 29413502
$ws.Range("K16").Value = 29413502
$ws.Range("M16").Value = -29413332
$ws.Range("H61").Value = 42133.883
$ws.Range("I61").Value = 45299.207
$ws.Range("K61").Value = 45299.207
$ws.Range("M61").Value = -45097.207
$ws.Range("H113").Value = 42133.883
$ws.Range("I113").Value = 45299.207
$ws.Range("K113").Value = 45299.207
$ws.Range("M113").Value = -43129.207
$ws.Range("H138").Value = 96095.25
$ws.Range("J138").Value = 96095.25
$ws.Range("L138").Value = 96095.25
$ws.Range("N138").Value = -106375.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9631.157999999999
$ws.Range("I81").Value = 18832
$ws.Range("K81").Value = 37664
$ws.Range("M81").Value = -36603
$ws.Range("H84").Value = 9631.157999999999
$ws.Range("I84").Value = 18832
$ws.Range("K84").Value = 188320
$ws.Range("M84").Value = -183016
$ws.Range("H97").Value = 34840
$ws.Range("I97").Value = 10520
$ws.Range("K97").Value = 10520
$ws.Range("M97").Value = -9529
